$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Responsable "Oriana Osiris" -> "Marisol Ornelas" for rows 4-7
$ws.Range("C4:C7").Value = "Marisol Ornelas"

# Status "En proceso" -> "Cerrada" for rows 4-6 (row 7 is already "Cerrada")
$ws.Range("F4:F6").Value = "Cerrada"

# Fecha real de cierre filled in for rows 4-7
$ws.Range("E4:E7").Value = 42383

# Row 8 - responsable cleared entirely
$ws.Range("C8").Value = ""

# Update selection to D7
$ws.Range("D7").Select()
